# "show wide data from db into second html table"
# Update the sample statistics row with the new accuracy percentages coming
# back from the (now wider) DB-backed second table, and leave the sheet
# scrolled/selected where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 43
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 57
$ws.Range("J2").Value = 37

# Scroll the view so column E is the left-most visible column, then leave the
# selection on J12 (matches the author's last on-screen position).
$ws.Range("J12").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
